$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.651.19"
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = "  +0.36%  "

# Row 3
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.963.19"
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = "  +2.34%  "

# Row 4
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.81"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  +1.65%  "

# Row 6
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4841"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  +0.74%  "

# Row 8
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2948"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +1.67%  "

# Row 9
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06790"
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = "  +1.01%  "

# Row 10
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "109.97"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -0.54%  "

# Row 11
$ws.Range("E11").Value = "  +1.79%  "

# Row 12
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.974.30"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  +2.97%  "

# Row 13
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07770"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  +2.62%  "

# Row 14
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.462"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +3.59%  "

# Row 15
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6887"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +3.09%  "

# Row 16
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "294.23"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  -1.39%  "

# Row 17
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.679.57"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +0.48%  "

# Row 18
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.25"
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = "  +1.88%  "

# Row 19
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.229.67"
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = "  +2.84%  "

# Row 20
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007706"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +1.69%  "

# Row 21
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.617"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +0.64%  "

# Row 22
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9993"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  +0.05%  "

# Row 24
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.605"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  +2.02%  "

# Row 25
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.889"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +4.20%  "

# Row 26
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.42"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +3.56%  "

# Row 27
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.22"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -0.42%  "

# Row 28
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.181"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +3.15%  "

# Row 29
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1065"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -1.17%  "

# Row 30
$ws.Range("E30").Value = "  +2.65%  "

# Row 31
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.718"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +16.64%  "

# Row 32
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.444"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +6.64%  "

# Row 33
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05117"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +2.34%  "

# Row 34
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7710"
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = "  +4.66%  "

# Row 35
$ws.Range("E35").Value = "  +3.86%  "

# Row 36
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02048"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +0.73%  "

# Row 37
$ws.Range("E37").Value = "  +0.36%  "

# Row 38
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.722"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +1.38%  "

# Row 39
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.123"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  +5.08%  "

# Row 40
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.386"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +8.28%  "

# Row 41
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4477"
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = "  +0.98%  "

# Row 42
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.08"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -1.76%  "

# Row 43
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8766"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  +1.70%  "

# Row 44
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.40"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -2.66%  "

# Row 45
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +0.02%  "

# Row 46
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.495"
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = "  +2.99%  "

# Row 47
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1282"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +4.19%  "

# Row 48
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.405"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +0.93%  "

# Row 49
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.04"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  +2.74%  "

# Row 50
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "47.65"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  -3.47%  "

# Row 51
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "917.47"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +6.32%  "

